$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$values = @(
    @("-60.452222850300473", "-46.815454639102917", "32.234359752781799", "17.574269392807071", "54.662865537799988", "23.739806067855504", "8.0253940776943547E-2", "3.9251599338186678E-2", "0.1792914238915696"),
    @("-38.568789654826901", "-37.058758765709939", "39.652384484547859", "22.233071697874685", "112.06024561874294", "44.848392915193649", "0.77604792271659129", "2.338784622668174E-2", "2.6059326658168607E-2"),
    @("-32.855435271891913", "-18.719792706390411", "24.886237802842782", "27.028819409305996", "73.584706403222611", "98.551787665001754", "1.5320075979370769E-2", "0.8090047428123075", "0.51605489471206889"),
    @("-97.92097308662197", "-74.541131221758448", "23.813821668604032", "24.80804778585167", "23.508923601223366", "33.700803808855369", "1.990728685446809E-2", "0.9249022513445011", "0.56478935286805643")
)

for ($r = 0; $r -lt 4; $r++) {
    for ($c = 0; $c -lt 9; $c++) {
        $ws1.Cells.Item($r + 1, $c + 1).Value = [double]$values[$r][$c]
    }
}

$excel.Calculate()
